# Weekly driver report update for 2025-04-28
# Insert a new "Good Driver" row at the top of the good-drivers table
# (row 12) for a newly-seen driver version, and bump the client-count
# (Total Samples) figures for the drivers that were already in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing good-driver rows (12-17) down one row to make room
# for the new driver entry.
$ws.Rows.Item(12).Insert()

# New driver row. This driver version has no recorded "Driver Vintage"
# date yet, so the date column is left blank.
$ws.Cells.Item(12, 1).Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Cells.Item(12, 2).Value = 11128
$ws.Cells.Item(12, 2).NumberFormat = "#,##0"
$ws.Cells.Item(12, 2).HorizontalAlignment = -4152
$ws.Cells.Item(12, 4).Value = 100
$ws.Cells.Item(12, 5).Value = 0

# Updated client counts ("Total Samples") for the pre-existing drivers,
# now shifted down to rows 13-18.
$ws.Cells.Item(13, 2).Value = 486214
$ws.Cells.Item(14, 2).Value = 79953
$ws.Cells.Item(15, 2).Value = 35355
$ws.Cells.Item(16, 2).Value = 65425
$ws.Cells.Item(17, 2).Value = 117653
